$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kosong")

$ws.Range("A22").Value = 6065
$ws.Range("A23").Value = 6066
$ws.Range("A24").Value = 6067
$ws.Range("A25").Value = 7796
$ws.Range("A26").Value = 7797
$ws.Range("A27").Value = 7798
$ws.Range("A28").Value = 9170
$ws.Range("A29").Value = 9171
